# Fruta / hortaliza, semanal
# Re-assigns the weekly records (Fecha, Volumen, Precio mínimo/máximo/promedio,
# Unidad de comercialización, Precio $/Kg) to different rows - the underlying
# set of weekly observations is unchanged, only which row holds which week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D = 44400; M = 100; N = 10000; O = 10000; P = 10000; Q = '$/caja 14 kilos';           S = 714 }
    3 = @{ D = 44208; M = 210; N = 10000; O = 10000; P = 10000; Q = '$/caja 14 kilos empedrada'; S = 714 }
    4 = @{ D = 44309; M = 300; N = 7000;  O = 7000;  P = 7000;  Q = '$/caja 14 kilos empedrada'; S = 500 }
    5 = @{ D = 44397; M = 60;  N = 11000; O = 11000; P = 11000; Q = '$/caja 14 kilos';           S = 786 }
    6 = @{ D = 44162; M = 120; N = 7000;  O = 7000;  P = 7000;  Q = '$/caja 14 kilos empedrada'; S = 500 }
    7 = @{ D = 44351; M = 300; N = 10000; O = 10000; P = 10000; Q = '$/caja 14 kilos empedrada'; S = 714 }
    8 = @{ D = 44176; M = 250; N = 7000;  O = 7000;  P = 7000;  Q = '$/caja 14 kilos empedrada'; S = 500 }
    9 = @{ D = 44491; M = 180; N = 9000;  O = 9000;  P = 9000;  Q = '$/caja 14 kilos empedrada'; S = 643 }
}

foreach ($row in $data.Keys) {
    $rec = $data[$row]
    $ws.Range("D$row").Value = $rec.D
    $ws.Range("M$row").Value = $rec.M
    $ws.Range("N$row").Value = $rec.N
    $ws.Range("O$row").Value = $rec.O
    $ws.Range("P$row").Value = $rec.P
    $ws.Range("Q$row").Value = $rec.Q
    $ws.Range("S$row").Value = $rec.S
}
